$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row to append, one past the current last row (68 -> 69).
$newRow = 69
$lastRow = $newRow - 1

# Copy the formatting (styles/borders/number formats) of the last existing
# data row so the appended row matches the sheet's look (bold/bordered
# index column, date-time formatted match-date column, etc.).
$ws.Range("A" + $lastRow + ":V" + $lastRow).Copy($ws.Range("A" + $newRow + ":V" + $newRow))

# Populate the new row with the match data.
$ws.Cells.Item($newRow, 1).Value  = 68
$ws.Cells.Item($newRow, 2).Value  = "croatia"
$ws.Cells.Item($newRow, 3).Value  = "hnl"
$ws.Cells.Item($newRow, 4).Value  = "2023-2024"
$ws.Cells.Item($newRow, 5).Value  = 45240.75
$ws.Cells.Item($newRow, 6).Value  = "Slaven Belupo"
$ws.Cells.Item($newRow, 7).Value  = 0
$ws.Cells.Item($newRow, 8).Value  = "Gorica"
$ws.Cells.Item($newRow, 9).Value  = 0
$ws.Cells.Item($newRow, 10).Value = 2.8
$ws.Cells.Item($newRow, 11).Value = "05/11/2023 17:12"
$ws.Cells.Item($newRow, 12).Value = 2.93
$ws.Cells.Item($newRow, 13).Value = "10/11/2023 17:59"
$ws.Cells.Item($newRow, 14).Value = 3.19
$ws.Cells.Item($newRow, 15).Value = "05/11/2023 17:12"
$ws.Cells.Item($newRow, 16).Value = 3.1
$ws.Cells.Item($newRow, 17).Value = "10/11/2023 17:35"
$ws.Cells.Item($newRow, 18).Value = 2.62
$ws.Cells.Item($newRow, 19).Value = "05/11/2023 17:12"
$ws.Cells.Item($newRow, 20).Value = 2.63
$ws.Cells.Item($newRow, 21).Value = "10/11/2023 17:59"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/croatia/hnl/slaven-belupo-hnk-gorica/IFF6CVmo/"

Write-Host ("Added row " + $newRow + " to sheet '" + $ws.Name + "'.")
